# Update column G ("K") values per regenerated save_data
# (save_data regen to use K instead of Strike#, regen std/mean, calc and write s_vals)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 3
    4  = 1
    5  = 5
    6  = 7
    7  = 3
    8  = 6
    9  = 3
    10 = 1
    11 = 7
    12 = 3
    13 = 6
    14 = 3
    15 = 4
    16 = 7
    17 = 2
    18 = 5
    19 = 4
    20 = 6
    21 = 5
    22 = 5
    23 = 1
}

foreach ($row in $values.Keys) {
    $ws.Range("G$row").Value = $values[$row]
}
